$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column A width (stored XML width of 13 == ColumnWidth 13 - 5/6 padding offset)
$ws.Columns.Item(1).ColumnWidth = 12.166666666666666

# New rows 6, 7, 9, 10
$ws.Range("A6").Value = "ADC Ref"
$ws.Range("B6").Value = 5

$ws.Range("A7").Value = "data"
$ws.Range("B7").Value = 112
$ws.Range("D7").Formula = "=(B7/1024 * B6) *B4"

$ws.Range("A9").Value = "Calculated Voltage 1"
$ws.Range("B9").Formula = "=FLOOR(B7*B4 / 1024,1) *B6"

$ws.Range("A10").Value = "Calculated Voltage 2"
$ws.Range("B10").Formula = "=FLOOR(B7*B4*B6 /1024,1)"

# Update selection to B7
$ws.Range("B7").Select()
